# Apply timetable updates (new course assignments) and add extra daily time slots
# (12:00-13:00, 13:00-14:00, 15:30-16:30, 16:30-17:30, 17:30-18:30) to the
# Section_A and Section_B sheets, reflecting the change of tutorial duration
# from 1.5hr blocks to 1hr blocks. Also restructure the Course_Summary sheet
# to separate "Credits" into "Lectures/Week", "Tutorials/Week" and
# "Total Credits" columns, and update instructor names.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("A2").Value = "9:00-10:30"
$ws.Range("B2").Value = "EC101"
$ws.Range("C2").Value = "EC101"
$ws.Range("D2").Value = "CS101"
$ws.Range("E2").Value = "EC101"
$ws.Range("F2").Value = "CS101"
$ws.Range("A3").Value = "10:30-12:00"
$ws.Range("B3").Value = "HS101"
$ws.Range("C3").Value = "Free"
$ws.Range("D3").Value = "Free"
$ws.Range("E3").Value = "MA101"
$ws.Range("F3").Value = "Free"
$ws.Range("A5").Value = "14:00-15:30"
$ws.Range("B5").Value = "Free"
$ws.Range("C5").Value = "MA102"
$ws.Range("D5").Value = "HS101"
$ws.Range("E5").Value = "Free"
$ws.Range("F5").Value = "Free"
$ws.Range("A6").Value = "15:30-17:00"
$ws.Range("B6").Value = "DS101"
$ws.Range("C6").Value = "DS101"
$ws.Range("D6").Value = "Free"
$ws.Range("E6").Value = "CS151 (Elective)"
$ws.Range("F6").Value = "Free"
$ws.Range("A7").Value = "17:00-18:30"
$ws.Range("B7").Value = "CS101"
$ws.Range("C7").Value = "MA101"
$ws.Range("D7").Value = "Free"
$ws.Range("E7").Value = "HS101"
$ws.Range("F7").Value = "MA102"
$ws.Range("A8").Value = "12:00-13:00"
$ws.Range("B8").Value = "Free"
$ws.Range("C8").Value = "Free"
$ws.Range("D8").Value = "Free"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "Free"
$ws.Range("A9").Value = "13:00-14:00"
$ws.Range("B9").Value = "Free"
$ws.Range("C9").Value = "Free"
$ws.Range("D9").Value = "Free"
$ws.Range("E9").Value = "Free"
$ws.Range("F9").Value = "Free"
$ws.Range("A10").Value = "15:30-16:30"
$ws.Range("B10").Value = "Free"
$ws.Range("C10").Value = "Free"
$ws.Range("D10").Value = "Free"
$ws.Range("E10").Value = "Free"
$ws.Range("F10").Value = "Free"
$ws.Range("A11").Value = "16:30-17:30"
$ws.Range("B11").Value = "Free"
$ws.Range("C11").Value = "Free"
$ws.Range("D11").Value = "Free"
$ws.Range("E11").Value = "Free"
$ws.Range("F11").Value = "Free"
$ws.Range("A12").Value = "17:30-18:30"
$ws.Range("B12").Value = "Free"
$ws.Range("C12").Value = "Free"
$ws.Range("D12").Value = "Free"
$ws.Range("E12").Value = "Free"
$ws.Range("F12").Value = "Free"
$newRowsA = $ws.Range("A8:A12")
$newRowsA.Font.Bold = $true
$newRowsA.Borders.LineStyle = 1
$newRowsA.HorizontalAlignment = -4108
$newRowsA.VerticalAlignment = -4160

$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("A2").Value = "9:00-10:30"
$ws.Range("B2").Value = "MA101"
$ws.Range("C2").Value = "CS101"
$ws.Range("D2").Value = "DS101"
$ws.Range("E2").Value = "Free"
$ws.Range("F2").Value = "MA102"
$ws.Range("A3").Value = "10:30-12:00"
$ws.Range("B3").Value = "Free"
$ws.Range("C3").Value = "Free"
$ws.Range("D3").Value = "Free"
$ws.Range("E3").Value = "Free"
$ws.Range("F3").Value = "EC101"
$ws.Range("A5").Value = "14:00-15:30"
$ws.Range("B5").Value = "Free"
$ws.Range("C5").Value = "EC101"
$ws.Range("D5").Value = "Free"
$ws.Range("E5").Value = "EC101"
$ws.Range("F5").Value = "CS101"
$ws.Range("A6").Value = "15:30-17:00"
$ws.Range("B6").Value = "MA102"
$ws.Range("C6").Value = "MA101"
$ws.Range("D6").Value = "HS101"
$ws.Range("E6").Value = "CS151 (Elective)"
$ws.Range("F6").Value = "DS101"
$ws.Range("A7").Value = "17:00-18:30"
$ws.Range("B7").Value = "HS101"
$ws.Range("C7").Value = "Free"
$ws.Range("D7").Value = "CS101"
$ws.Range("E7").Value = "Free"
$ws.Range("F7").Value = "HS101"
$ws.Range("A8").Value = "12:00-13:00"
$ws.Range("B8").Value = "Free"
$ws.Range("C8").Value = "Free"
$ws.Range("D8").Value = "Free"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "Free"
$ws.Range("A9").Value = "13:00-14:00"
$ws.Range("B9").Value = "Free"
$ws.Range("C9").Value = "Free"
$ws.Range("D9").Value = "Free"
$ws.Range("E9").Value = "Free"
$ws.Range("F9").Value = "Free"
$ws.Range("A10").Value = "15:30-16:30"
$ws.Range("B10").Value = "Free"
$ws.Range("C10").Value = "Free"
$ws.Range("D10").Value = "Free"
$ws.Range("E10").Value = "Free"
$ws.Range("F10").Value = "Free"
$ws.Range("A11").Value = "16:30-17:30"
$ws.Range("B11").Value = "Free"
$ws.Range("C11").Value = "Free"
$ws.Range("D11").Value = "Free"
$ws.Range("E11").Value = "Free"
$ws.Range("F11").Value = "Free"
$ws.Range("A12").Value = "17:30-18:30"
$ws.Range("B12").Value = "Free"
$ws.Range("C12").Value = "Free"
$ws.Range("D12").Value = "Free"
$ws.Range("E12").Value = "Free"
$ws.Range("F12").Value = "Free"
$newRowsA = $ws.Range("A8:A12")
$newRowsA.Font.Bold = $true
$newRowsA.Borders.LineStyle = 1
$newRowsA.HorizontalAlignment = -4108
$newRowsA.VerticalAlignment = -4160

$ws = $wb.Worksheets.Item("Course_Summary")
$oldInstructorHeader = $ws.Range("F1").Value()
$ws.Range("E1").Value = "Lectures/Week"
$ws.Range("F1").Value = "Tutorials/Week"
$ws.Range("G1").Value = "Total Credits"
$ws.Range("H1").Value = $oldInstructorHeader
$hdrRng = $ws.Range("G1:H1")
$hdrRng.Font.Bold = $true
$hdrRng.Borders.LineStyle = 1
$hdrRng.HorizontalAlignment = -4108
$hdrRng.VerticalAlignment = -4160
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = "Dr. Meera Nair"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Dr. Arjun Deshmukh"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "Dr. Chintamani"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = "Dr. Kavita Bansal"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = "Dr. Suresh Kulkarni"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = "Dr. Rajesh N S"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = "Dr. Neel Patel"
